# Weekly price update: a new week of "Albahaca" price data is published.
# The whole data table (rows 8-64) is shifted down by one row (the oldest
# rows move further down), a fresh row is inserted at row 8 for the newest
# report date (Fecha = 44532), and the table grows from A1:R64 to A1:R65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 8..64 down to 9..65 (this also duplicates the trailing row,
# extending the used range to row 65 and bumping the sheet dimension).
$ws.Rows("8:8").Insert()

# Populate the newly-opened row 8 with this week's record.
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44532
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112052
$ws.Range("G8").Value = "Albahaca"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 740
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 3500
$ws.Range("N8").Value = "$/paquete"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 3500
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"
